$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to text format so numeric-looking strings
# (e.g. "1.00", "3.722.29", "  +0.12%  ") are preserved exactly as text,
# matching the original inline-string cell contents.
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = "68.966.94"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "3.770.71"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "629.84"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("D6").Value = "165.40"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").Value = "3.767.44"
$ws.Range("E7").Value = "  -1.35%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("E10").Value = "  -1.64%  "
$ws.Range("D11").Value = "0.457"
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").Value = "6.77"
$ws.Range("E12").Value = "  +1.85%  "
$ws.Range("E13").Value = "  -4.30%  "
$ws.Range("D14").Value = "34.80"
$ws.Range("E14").Value = "  -3.01%  "
$ws.Range("D15").Value = "4.403.08"
$ws.Range("E15").Value = "  -1.23%  "
$ws.Range("D16").Value = "3.782.94"
$ws.Range("E16").Value = "  -2.18%  "
$ws.Range("D17").Value = "68.964.41"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").Value = "17.66"
$ws.Range("E18").Value = "  -2.70%  "
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("E20").Value = "  -2.02%  "
$ws.Range("D21").Value = "467.42"
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("D22").Value = "9.52"
$ws.Range("E22").Value = "  -1.69%  "
$ws.Range("D23").Value = "0.701"
$ws.Range("E23").Value = "  -0.64%  "
$ws.Range("D24").Value = "82.03"
$ws.Range("E24").Value = "  -2.17%  "
$ws.Range("E25").Value = "  -6.76%  "
$ws.Range("D26").Value = "12.10"
$ws.Range("E26").Value = "  +0.72%  "
$ws.Range("E27").Value = "  -1.44%  "
$ws.Range("E28").Value = "  +0.61%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").Value = "3.920.15"
$ws.Range("E30").Value = "  -1.25%  "
$ws.Range("D31").Value = "2.67"
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("E32").Value = "  +0.74%  "
$ws.Range("E33").Value = "  -3.47%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "28.40"
$ws.Range("E34").Value = "  -2.57%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "0.175"
$ws.Range("E35").Value = "  +18.15%  "
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("D37").Value = "3.722.29"
$ws.Range("E37").Value = "  -1.14%  "
$ws.Range("E38").Value = "  -2.42%  "
$ws.Range("E39").Value = "  -0.96%  "
$ws.Range("D40").Value = "3.26"
$ws.Range("E40").Value = "  -2.15%  "
$ws.Range("E41").Value = "  -2.44%  "
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").Value = "0.963"
$ws.Range("E42").Value = "  -1.57%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("E45").Value = "  +5.13%  "
$ws.Range("D46").Value = "155.83"
$ws.Range("E46").Value = "  +0.74%  "
$ws.Range("D47").Value = "43.37"
$ws.Range("E47").Value = "  +1.09%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").Value = "46.93"
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").Value = "1.40"
$ws.Range("E49").Value = "  -1.83%  "
$ws.Range("E50").Value = "  -2.59%  "
$ws.Range("E51").Value = "  -1.30%  "

# Restore default cell style (no explicit style index) on the Price/Volume range
# now that values have been written, so formatting matches the source workbook.
$priceVolRange.Style = "Normal"
